$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: "12 ماهه منتهی به ..." period labels, rolled forward one year ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/10"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/10"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/10"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/10"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/10"

# --- Header row 9: "تاریخ انتشار" publish-date labels, rolled forward one year ---
$ws.Range("D9").Value = "1399-02-31 (12)"
$ws.Range("E9").Value = "1400-02-27 (9)"
$ws.Range("F9").Value = "1401-02-25 (8)"
$ws.Range("G9").Value = "1402-02-24 (10)"
$ws.Range("H9").Value = "1402-02-30 (3)"

# --- Row 11: فروش ---
$ws.Range("D11").Value = 11605
$ws.Range("E11").Value = 11934
$ws.Range("F11").Value = 8507
$ws.Range("G11").Value = 13640
$ws.Range("H11").Value = 19856

# --- Row 12 ---
$ws.Range("D12").Value = -8937
$ws.Range("E12").Value = -9418
$ws.Range("F12").Value = -6642
$ws.Range("G12").Value = -9314
$ws.Range("H12").Value = -12733

# --- Row 13 ---
$ws.Range("D13").Value = 2668
$ws.Range("E13").Value = 2516
$ws.Range("F13").Value = 1865
$ws.Range("G13").Value = 4326
$ws.Range("H13").Value = 7123

# --- Row 14 ---
$ws.Range("D14").Value = -1679
$ws.Range("E14").Value = -1992
$ws.Range("F14").Value = -898
$ws.Range("G14").Value = -1478
$ws.Range("H14").Value = -2712

# --- Row 15 ---
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = -446
$ws.Range("G15").Value = -309
$ws.Range("H15").Value = "-"

# --- Row 16 ---
$ws.Range("D16").Value = 55
$ws.Range("E16").Value = -1060
$ws.Range("F16").Value = -199
$ws.Range("G16").Value = -54
$ws.Range("H16").Value = -186

# --- Row 17 ---
$ws.Range("D17").Value = 1044
$ws.Range("E17").Value = -536
$ws.Range("F17").Value = 322
$ws.Range("G17").Value = 2486
$ws.Range("H17").Value = 4225

# --- Row 18 ---
$ws.Range("D18").Value = -1055
$ws.Range("E18").Value = -799
$ws.Range("F18").Value = -373
$ws.Range("G18").Value = -630
$ws.Range("H18").Value = -126

# --- Row 19 ---
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = 293
$ws.Range("F19").Value = 238
$ws.Range("G19").Value = -26
$ws.Range("H19").Value = 492

# --- Row 20 ---
$ws.Range("D20").Value = 28
$ws.Range("E20").Value = -1043
$ws.Range("F20").Value = 187
$ws.Range("G20").Value = 1830
$ws.Range("H20").Value = 4591

# --- Row 21 ---
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = -76
$ws.Range("G21").Value = -591
$ws.Range("H21").Value = -939

# --- Row 22 ---
$ws.Range("D22").Value = 28
$ws.Range("E22").Value = -1043
$ws.Range("F22").Value = 111
$ws.Range("G22").Value = 1239
$ws.Range("H22").Value = 3653

# --- Row 23 is unchanged (all "-") ---

# --- Row 24 ---
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -1043
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 1239
$ws.Range("H24").Value = 3653

# --- Row 25 is unchanged (all 0) ---

# --- Row 26 ---
$ws.Range("D26").Value = 4634
$ws.Range("E26").Value = 3294
$ws.Range("F26").Value = 5728
$ws.Range("G26").Value = 4625
$ws.Range("H26").Value = 3782

# --- Row 27 is unchanged (all 0) ---
